$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 710
$ws.Range("I58").Value = 710
$ws.Range("K58").Value = 2130
$ws.Range("M58").Value = -1980
$ws.Range("H74").Value = 3753.2942
$ws.Range("I74").Value = 3311.7778
$ws.Range("J74").Value = 4250
$ws.Range("K74").Value = 3311.7778
$ws.Range("L74").Value = 4250
$ws.Range("M74").Value = -2375.7778
$ws.Range("N74").Value = -6122
$ws.Range("H76").Value = 3441.4666
$ws.Range("I76").Value = 2973.5278
$ws.Range("J76").Value = 5313.222
$ws.Range("K76").Value = 2973.5278
$ws.Range("L76").Value = 5313.222
$ws.Range("M76").Value = -2658.5278
$ws.Range("N76").Value = -5943.222
$ws.Range("H77").Value = 3753.2942
$ws.Range("I77").Value = 3311.7778
$ws.Range("J77").Value = 4250
$ws.Range("K77").Value = 16558.889
$ws.Range("L77").Value = 21250
$ws.Range("M77").Value = -11878.889
$ws.Range("N77").Value = -30610
$ws.Range("H79").Value = 3441.4666
$ws.Range("I79").Value = 2973.5278
$ws.Range("J79").Value = 5313.222
$ws.Range("K79").Value = 2973.5278
$ws.Range("L79").Value = 5313.222
$ws.Range("M79").Value = -1881.5278
$ws.Range("N79").Value = -7497.222
$ws.Range("H86").Value = 126276.5
$ws.Range("J86").Value = 201100
$ws.Range("L86").Value = 201100
$ws.Range("N86").Value = -203346
$ws.Range("H89").Value = 126276.5
$ws.Range("J89").Value = 201100
$ws.Range("L89").Value = 1005500
$ws.Range("N89").Value = -1016732
$ws.Range("H96").Value = 50003030
$ws.Range("I96").Value = 62502532
$ws.Range("J96").Value = 5029
$ws.Range("K96").Value = 187507596
$ws.Range("L96").Value = 15087
$ws.Range("M96").Value = -187506223
$ws.Range("N96").Value = -17833
$ws.Range("H98").Value = 834.6818
$ws.Range("I98").Value = 834.6818
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 834.6818
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 663.3182
$ws.Range("H105").Value = 43557
$ws.Range("J105").Value = 43557
$ws.Range("L105").Value = 43557
$ws.Range("N105").Value = -50545
$ws.Range("H112").Value = 1119.421
$ws.Range("J112").Value = 1265.2667
$ws.Range("L112").Value = 3795.800099999999
$ws.Range("N112").Value = -6011.800099999999
$ws.Range("H122").Value = 834.6818
$ws.Range("I122").Value = 834.6818
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2504.0454
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -54.04539999999997
$ws.Range("H125").Value = 906446.9
$ws.Range("I125").Value = 1972
$ws.Range("J125").Value = 1471743.6
$ws.Range("K125").Value = 17748
$ws.Range("L125").Value = 13245692.4
$ws.Range("M125").Value = -15288
$ws.Range("N125").Value = -13250612.4
$ws.Range("H138").Value = 4951.5415
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4951.5415
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 14854.6245
$ws.Range("N138").Value = -25134.6245
$ws.Range("N98").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 2000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1314
$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 2000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 10000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6568
$ws.Range("H88").Value = 1750
$ws.Range("I88").Value = 1500
$ws.Range("K88").Value = 1500
$ws.Range("M88").Value = -1094
$ws.Range("H91").Value = 1750
$ws.Range("I91").Value = 1500
$ws.Range("K91").Value = 1500
$ws.Range("M91").Value = -96
$ws.Range("H97").Value = 2750.8667
$ws.Range("I97").Value = 2096.3076
$ws.Range("J97").Value = 7005.5
$ws.Range("K97").Value = 2096.3076
$ws.Range("L97").Value = 7005.5
$ws.Range("M97").Value = -1600.3076
$ws.Range("N97").Value = -7997.5
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1696.091
$ws.Range("I86").Value = 1761.8889
$ws.Range("K86").Value = 1761.8889
$ws.Range("M86").Value = -638.8888999999999
$ws.Range("H89").Value = 1696.091
$ws.Range("I89").Value = 1761.8889
$ws.Range("K89").Value = 8809.4445
$ws.Range("M89").Value = -3193.4445
$ws.Range("H94").Value = 1142.9
$ws.Range("I94").Value = 1118.625
$ws.Range("J94").Value = 1240
$ws.Range("K94").Value = 1118.625
$ws.Range("L94").Value = 1240
$ws.Range("M94").Value = -667.625
$ws.Range("N94").Value = -2142
$ws.Range("H99").Value = 1772.6086
$ws.Range("I99").Value = 1224.5454
$ws.Range("K99").Value = 1224.5454
$ws.Range("M99").Value = 273.4546
$ws.Range("H105").Value = 5152.5
$ws.Range("I105").Value = 5152.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5152.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3405.5
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4173.8667
$ws.Range("I62").Value = 2900.5
$ws.Range("J62").Value = 5629.143
$ws.Range("K62").Value = 2900.5
$ws.Range("L62").Value = 5629.143
$ws.Range("M62").Value = -2276.5
$ws.Range("N62").Value = -6877.143
$ws.Range("H65").Value = 4173.8667
$ws.Range("I65").Value = 2900.5
$ws.Range("J65").Value = 5629.143
$ws.Range("K65").Value = 14502.5
$ws.Range("L65").Value = 28145.715
$ws.Range("M65").Value = -11382.5
$ws.Range("N65").Value = -34385.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 791.8182
$ws.Range("I46").Value = 150
$ws.Range("J46").Value = 934.44446
$ws.Range("K46").Value = 450
$ws.Range("L46").Value = 2803.33338
$ws.Range("M46").Value = -359
$ws.Range("N46").Value = -2985.33338
$ws.Range("H133").Value = 1977.7858
$ws.Range("I133").Value = 1961.125
$ws.Range("K133").Value = 5883.375
$ws.Range("M133").Value = -823.375
$ws.Range("H136").Value = 1954.8462
$ws.Range("I136").Value = 1867.4166
$ws.Range("J136").Value = 1993.7037
$ws.Range("K136").Value = 5602.2498
$ws.Range("L136").Value = 5981.1111
$ws.Range("M136").Value = -502.2497999999996
$ws.Range("N136").Value = -16181.1111
$ws.Range("H138").Value = 1474.5834
$ws.Range("I138").Value = 1030.7693
$ws.Range("K138").Value = 3092.3079
$ws.Range("M138").Value = 2047.6921
$ws.Range("H139").Value = 918.2222
$ws.Range("I139").Value = 825.17645
$ws.Range("J139").Value = 2500
$ws.Range("K139").Value = 2475.52935
$ws.Range("L139").Value = 7500
$ws.Range("M139").Value = 2664.47065
$ws.Range("N139").Value = -17780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2405.5557
$ws.Range("I97").Value = 3252
$ws.Range("J97").Value = 1347.5
$ws.Range("K97").Value = 3252
$ws.Range("L97").Value = 1347.5
$ws.Range("M97").Value = -2756
$ws.Range("N97").Value = -2339.5
$ws.Range("H122").Value = 2766.2964
$ws.Range("I122").Value = 2562.5557
$ws.Range("J122").Value = 3173.7778
$ws.Range("K122").Value = 7687.6671
$ws.Range("L122").Value = 9521.3334
$ws.Range("M122").Value = -5237.6671
$ws.Range("N122").Value = -14421.3334
$ws.Range("H136").Value = 19400
$ws.Range("J136").Value = 19400
$ws.Range("L136").Value = 58200
$ws.Range("N136").Value = -63300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 11008.182
$ws.Range("J93").Value = 2450.75
$ws.Range("L93").Value = 2450.75
$ws.Range("N93").Value = -4946.75
$ws.Range("H132").Value = 7951.048
$ws.Range("I132").Value = 9063.925999999999
$ws.Range("J132").Value = 5947.8667
$ws.Range("K132").Value = 27191.778
$ws.Range("L132").Value = 17843.6001
$ws.Range("M132").Value = -24661.778
$ws.Range("N132").Value = -22903.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1400
$ws.Range("I96").Value = 1050
$ws.Range("J96").Value = 1750
$ws.Range("K96").Value = 1050
$ws.Range("L96").Value = 1750
$ws.Range("M96").Value = 323
$ws.Range("N96").Value = -4496
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("H137").Value = 80600
$ws.Range("J137").Value = 80600
$ws.Range("L137").Value = 80600
$ws.Range("N137").Value = -90800
$ws.Range("N97").ClearContents()
